# This script applies the latest crypto price/volume snapshot to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.153.34"
Set-TextValue "E2" "  -0.35%  "

Set-TextValue "D3" "1.623.71"
Set-TextValue "E3" "  -1.32%  "

Set-TextValue "E4" "  -0.02%  "

Set-TextValue "D5" "214.33"
Set-TextValue "E5" "  -1.33%  "

Set-TextValue "E6" "  +1.43%  "

Set-TextValue "E7" "  -0.01%  "

Set-TextValue "E8" "  -1.33%  "

Set-TextValue "E9" "  -0.12%  "

Set-TextValue "D10" "20.33"
Set-TextValue "E10" "  +1.31%  "

Set-TextValue "E11" "  -0.10%  "

Set-TextValue "D12" "1.622.71"
Set-TextValue "E12" "  -1.28%  "

Set-TextValue "E13" "  -0.44%  "

Set-TextValue "E14" "  -0.58%  "

Set-TextValue "D15" "27.130.55"
Set-TextValue "E15" "  -0.37%  "

Set-TextValue "D16" "64.60"
Set-TextValue "E16" "  -4.18%  "

Set-TextValue "D17" "0.0₃0745"
Set-TextValue "E17" "  +0.40%  "

Set-TextValue "D18" "216.03"
Set-TextValue "E18" "  -1.75%  "

Set-TextValue "E19" "  -0.02%  "

Set-TextValue "E20" "  +0.59%  "

Set-TextValue "E21" "  -0.95%  "

Set-TextValue "E22" "  -6.52%  "

Set-TextValue "D23" "9.06"
Set-TextValue "E23" "  -1.93%  "

Set-TextValue "D24" "148.14"
Set-TextValue "E24" "  +0.26%  "

Set-TextValue "E25" "  -0.08%  "

Set-TextValue "D26" "7.29"
Set-TextValue "E26" "  -3.13%  "

Set-TextValue "E27" "  -1.31%  "

Set-TextValue "D28" "15.62"
Set-TextValue "E28" "  -1.20%  "

Set-TextValue "D29" "0.0508"
Set-TextValue "E29" "  -0.61%  "

Set-TextValue "E30" "  -1.11%  "

Set-TextValue "E31" "  -0.42%  "

Set-TextValue "E32" "  -1.16%  "

Set-TextValue "D33" "1.344.32"
Set-TextValue "E33" "  +4.54%  "

Set-TextValue "E34" "  -0.63%  "

Set-TextValue "E35" "  -0.64%  "

Set-TextValue "E36" "  +0.09%  "

Set-TextValue "E37" "  +1.63%  "

Set-TextValue "D38" "0.860"
Set-TextValue "E38" "  -0.37%  "

Set-TextValue "E39" "  -0.09%  "

Set-TextValue "E40" "  -0.64%  "

Set-TextValue "D41" "65.55"
Set-TextValue "E41" "  +5.85%  "

Set-TextValue "E42" "  -0.24%  "

Set-TextValue "E43" "  -1.38%  "

Set-TextValue "B44" "RocketPoolETH"
Set-TextValue "C44" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D44" "1.761.31"
Set-TextValue "E44" "  -1.37%  "

Set-TextValue "B45" "WEMIXToken"
Set-TextValue "C45" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D45" "0.926"
Set-TextValue "E45" "  +37.90%  "

Set-TextValue "D46" "90.05"
Set-TextValue "E46" "  -2.09%  "

Set-TextValue "E47" "  +1.16%  "

Set-TextValue "E48" "  +3.25%  "

Set-TextValue "E49" "  -0.32%  "

Set-TextValue "E50" "  -1.24%  "

Set-TextValue "E51" "  -0.04%  "
